# Apply the "add 2022-Q3 data" edit to the COP-康菲石油 workbook.
#
# 1. The "总计" (totals) summary sheet gets a new row for 2022-Q3 inserted
#    right after the header, pushing all the existing quarters down by one
#    row and adding a (new) trailing row for 2020-Q4.
# 2. A brand new detail worksheet named "2022-Q3" is inserted right after
#    "总计" (i.e. before "2022-Q2"), holding the per-fund breakdown for the
#    new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" summary sheet.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Data rows, top (most recent) to bottom, as they should look *after* the
# edit. Column A is a simple 0-based running index.
$totalsRows = @(
    @{ Label = "2022-Q3"; Count = 8;  Value = 3.13 },
    @{ Label = "2022-Q2"; Count = 6;  Value = 9.36 },
    @{ Label = "2022-Q1"; Count = 5;  Value = 6.03 },
    @{ Label = "2021-Q4"; Count = 5;  Value = 1.85 },
    @{ Label = "2021-Q3"; Count = 11; Value = 3.77 },
    @{ Label = "2021-Q2"; Count = 7;  Value = 4.18 },
    @{ Label = "2021-Q1"; Count = 4;  Value = 2.92 },
    @{ Label = "2020-Q4"; Count = 4;  Value = 1.71 }
)

# 2022-Q2's old D value (9.36) keeps its original floating point noise.
$totals.Range("D3").Value = 9.359999999999999

for ($i = 0; $i -lt $totalsRows.Count; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]

    # Column A (index) keeps the header/style formatting already present;
    # copy it down from row 8 (the last pre-existing row) the first time we
    # touch a brand new row (row 9).
    if ($r -eq 9) {
        $totals.Range("A8").Copy($totals.Range("A9"))
    }

    $totals.Range("A$r").Value = $i
    $totals.Range("B$r").Value = $row.Label
    $totals.Range("C$r").Value = $row.Count
    if (-not ($r -eq 3)) {
        $totals.Range("D$r").Value = $row.Value
    }
}

# ---------------------------------------------------------------------
# Step 2: create the new "2022-Q3" detail sheet (copy "2022-Q2" to inherit
# formatting/styling, then overwrite its contents).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Extend formatting for the two extra data rows (the copied sheet only has
# 6 data rows, the new quarter needs 8) by copying the last data row's
# look down twice.
$q3.Range("A7:H7").Copy($q3.Range("A8:H8"))
$q3.Range("A7:H7").Copy($q3.Range("A9:H9"))

$q3rows = @(
    @{ Code = "006679"; Name = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 A"; Size = "11.73";  Pos = "93.96"; Pct = "17.14"; Mv = "2.0105";  Rank = 1 },
    @{ Code = "162719"; Name = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A";        Size = "11.73";  Pos = "93.96"; Pct = "17.14"; Mv = "2.0105";  Rank = 1 },
    @{ Code = "006680"; Name = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 C"; Size = "5.92";   Pos = "93.96"; Pct = "17.14"; Mv = "1.0147";  Rank = 1 },
    @{ Code = "160416"; Name = "华安标普全球石油指数（QDII-LOF）A";                    Size = "2.74";   Pos = "93.58"; Pct = "3.64";  Mv = "0.0997";  Rank = 5 },
    @{ Code = "014982"; Name = "华安标普全球石油指数（QDII-LOF）C";                    Size = "0.22";   Pos = "93.58"; Pct = "3.64";  Mv = "0.0080";  Rank = 5 },
    @{ Code = "000049"; Name = "中银标普全球精选自然资源等权重指数（QDII）A";            Size = "0.23";   Pos = "88.52"; Pct = "0.97";  Mv = "0.0022";  Rank = 10 },
    @{ Code = "016334"; Name = "中银标普全球精选自然资源等权重指数（QDII）C";            Size = "0.00";   Pos = "88.52"; Pct = "0.97";  Mv = $null;      Rank = 10 },
    @{ Code = "004243"; Name = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C";         Size = "-11.74"; Pos = "93.96"; Pct = "17.14"; Mv = "-2.0122"; Rank = 1 }
)

# Fund code / name / size / position / pct / market-value columns hold
# text-formatted values (not real numbers) in the source data, so force
# the "@" (text) number format before writing them.
$q3.Range("B2:F9").NumberFormat = "@"

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $r = $i + 2
    $row = $q3rows[$i]

    $q3.Range("A$r").Value = $i
    $q3.Range("B$r").Value = $row.Code
    $q3.Range("C$r").Value = $row.Name
    $q3.Range("D$r").Value = $row.Size
    $q3.Range("E$r").Value = $row.Pos
    $q3.Range("F$r").Value = $row.Pct
    if ($row.Mv -eq $null) {
        $q3.Range("G$r").NumberFormat = "General"
        $q3.Range("G$r").Value = 0
    } else {
        $q3.Range("G$r").Value = $row.Mv
    }
    $q3.Range("H$r").Value = $row.Rank
}
